$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# Numeric value updates
$ws.Range("S6").Value = 1664
$ws.Range("S32").Value = 28307

# Text value updates (node names -> plain numeric-looking ids)
Set-TextValue "B8" "146"
Set-TextValue "B9" "78"
Set-TextValue "B10" "149"
Set-TextValue "B11" "82"
Set-TextValue "B12" "152"
Set-TextValue "B13" "86"
Set-TextValue "B14" "155"
Set-TextValue "B15" "90"
Set-TextValue "B16" "158"
Set-TextValue "B17" "94"
Set-TextValue "B18" "161"
Set-TextValue "B19" "121"
Set-TextValue "B20" "164"
Set-TextValue "B21" "167"
Set-TextValue "B22" "170"
Set-TextValue "B24" "173"
Set-TextValue "B25" "140"
Set-TextValue "B26" "176"

# Restore page margins to Excel defaults (inches -> points: 1pt = 1/72in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
